# Updates cryptos list: prices (D) and 1h volume % change (E) columns,
# plus two row swaps (Polygon/Uniswap at 22-23, Hedera/OKB at 38-39).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.650.46"
$ws.Range("E2").Value = "  +1.91%  "

$ws.Range("D3").Value = "'3.544.13"
$ws.Range("E3").Value = "  +0.73%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'608.58"
$ws.Range("E5").Value = "  +4.42%  "

$ws.Range("D6").Value = "'173.89"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7").Value = "'0.617"
$ws.Range("E7").Value = "  -0.72%  "

$ws.Range("D8").Value = "'3.541.76"
$ws.Range("E8").Value = "  +0.84%  "

$ws.Range("D10").Value = "'0.200"
$ws.Range("E10").Value = "  +5.08%  "

$ws.Range("D11").Value = "'6.75"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("D12").Value = "'0.586"
$ws.Range("E12").Value = "  -1.73%  "

$ws.Range("D13").Value = "'47.56"
$ws.Range("E13").Value = "  +0.98%  "

$ws.Range("D14").Value = "'0.0000280"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("D15").Value = "'4.113.92"
$ws.Range("E15").Value = "  +0.75%  "

$ws.Range("D16").Value = "'629.89"
$ws.Range("E16").Value = "  -6.92%  "

$ws.Range("D17").Value = "'8.46"
$ws.Range("E17").Value = "  -3.38%  "

$ws.Range("D18").Value = "'70.708.08"
$ws.Range("E18").Value = "  +1.93%  "

$ws.Range("D19").Value = "'3.548.44"
$ws.Range("E19").Value = "  +0.80%  "

$ws.Range("E20").Value = "  -1.83%  "

$ws.Range("D21").Value = "'17.45"
$ws.Range("E21").Value = "  -0.34%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'10.04"
$ws.Range("E22").Value = "  -10.55%  "

$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "'0.889"
$ws.Range("E23").Value = "  -1.79%  "

$ws.Range("D24").Value = "'15.89"
$ws.Range("E24").Value = "  -2.13%  "

$ws.Range("D25").Value = "'97.10"
$ws.Range("E25").Value = "  -1.09%  "

$ws.Range("E26").Value = "  -0.56%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("E28").Value = "  -1.95%  "

$ws.Range("D29").Value = "'9.21"
$ws.Range("E29").Value = "  -2.70%  "

$ws.Range("D30").Value = "'33.38"
$ws.Range("E30").Value = "  +0.96%  "

$ws.Range("D31").Value = "'3.12"
$ws.Range("E31").Value = "  -2.55%  "

$ws.Range("D32").Value = "'8.47"
$ws.Range("E32").Value = "  -3.38%  "

$ws.Range("E33").Value = "  -1.36%  "

$ws.Range("D34").Value = "'7.04"
$ws.Range("E34").Value = "  -4.50%  "

$ws.Range("D35").Value = "'566.91"
$ws.Range("E35").Value = "  -5.27%  "

$ws.Range("D36").Value = "'3.64"
$ws.Range("E36").Value = "  +0.56%  "

$ws.Range("D37").Value = "'10.79"
$ws.Range("E37").Value = "  -1.06%  "

$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'57.49"
$ws.Range("E38").Value = "  +0.30%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.102"
$ws.Range("E39").Value = "  -2.05%  "

$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("E41").Value = "  +5.65%  "

$ws.Range("E42").Value = "  +3.17%  "

$ws.Range("D43").Value = "'0.329"
$ws.Range("E43").Value = "  -2.48%  "

$ws.Range("D44").Value = "'3.335.85"
$ws.Range("E44").Value = "  -2.62%  "

$ws.Range("D45").Value = "'3.06"
$ws.Range("E45").Value = "  +4.66%  "

$ws.Range("D46").Value = "'0.0₃0717"
$ws.Range("E46").Value = "  +0.96%  "

$ws.Range("D47").Value = "'33.10"
$ws.Range("E47").Value = "  -1.14%  "

$ws.Range("D48").Value = "'2.65"
$ws.Range("E48").Value = "  +1.46%  "

$ws.Range("E49").Value = "  -2.58%  "

$ws.Range("D50").Value = "'134.18"
$ws.Range("E50").Value = "  +0.79%  "

$ws.Range("D51").Value = "'5.73"
$ws.Range("E51").Value = "  -1.34%  "
